$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("G5").Value = 2.38
$ws.Range("I5").Value = 3.6
$ws.Range("J5").Value = 3.25
$ws.Range("M5").Value = 1.14
$ws.Range("N5").Value = 5.5
$ws.Range("Q5").Value = 3.1
$ws.Range("R5").Value = 1.36
$ws.Range("Z5").Value = 9.5
$ws.Range("AF5").Value = 6
$ws.Range("AJ5").Value = 7.5
$ws.Range("AL5").Value = 15
$ws.Range("AN5").Value = 41
$ws.Range("AP5").Value = 2.14
$ws.Range("AQ5").Value = 1.68

# Row 6
$ws.Range("I6").Value = 4
$ws.Range("M6").Value = 1.11
$ws.Range("N6").Value = 6.5
$ws.Range("O6").Value = 1.5
$ws.Range("P6").Value = 2.5
$ws.Range("Q6").Value = 2.6
$ws.Range("R6").Value = 1.48
$ws.Range("U6").Value = 1.57
$ws.Range("V6").Value = 2.25
$ws.Range("AP6").Value = 2
$ws.Range("AQ6").Value = 1.85
$ws.Range("AR6").Value = 4.4
$ws.Range("AS6").Value = 1.2

# Row 7
$ws.Range("M7").Value = 1.1
$ws.Range("N7").Value = 7
$ws.Range("Q7").Value = 2.6
$ws.Range("R7").Value = 1.48
$ws.Range("S7").Value = 5.5
$ws.Range("T7").Value = 1.14
$ws.Range("AP7").Value = 1.98
$ws.Range("AQ7").Value = 1.88
$ws.Range("AR7").Value = 3.9
$ws.Range("AS7").Value = 1.24
